$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# This edit swaps the report rows for the "1d0d76ba..." and "d9605128..."
# e2e files (the d9605128 file moves up to row 5, 1d0d76ba moves to row 6)
# on all three worksheets, and updates the "d9605128" file's status from
# "Ready for handoff" to "In Translation" to reflect that it is now back
# in translation (archive/report regeneration).
# -------------------------------------------------------------------------

# ===================== Sheet 1: Overview =====================
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A5").Value = "d9605128-76f0-495d-8829-707337481735.md"
$ws1.Range("B5").Value = "e2e\d9605128-76f0-495d-8829-707337481735.md"
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-16 20:38:11"

$ws1.Range("A6").Value = "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md"
$ws1.Range("B6").Value = "e2e\1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md"
$ws1.Range("E6").Value = "Ready for handoff"
$ws1.Range("F6").Value = "Ready for handoff"
$ws1.Range("G6").Value = "2016-08-16 20:38:33"

$h1 = $ws1.Hyperlinks
$h1.Delete()
$h1.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b9dc80479d3c1d79b79aaf5b2b4c141c2df2ce0/e2e/4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md", [Type]::Missing, [Type]::Missing, "e2e\4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md")
$h1.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/49145225-09c9-46ae-803f-739b042b9bda.md", [Type]::Missing, [Type]::Missing, "e2e\49145225-09c9-46ae-803f-739b042b9bda.md")
$h1.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/b1b21d64-4f4b-4615-91cd-ab12a367e019.md", [Type]::Missing, [Type]::Missing, "e2e\b1b21d64-4f4b-4615-91cd-ab12a367e019.md")
$h1.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d6f18f4361b4819f466264699ed58e43a3b7325/e2e/d9605128-76f0-495d-8829-707337481735.md", [Type]::Missing, [Type]::Missing, "e2e\d9605128-76f0-495d-8829-707337481735.md")
$h1.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5995be60b2f8ebe30d1c30a9afe1cb90e32440ed/e2e/1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md", [Type]::Missing, [Type]::Missing, "e2e\1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md")
$h1.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d54b2f816603fc8dda517ffcd058267ebc1ea66/e2e/e12cf65d-4c4c-4743-92ec-438a640722aa.md", [Type]::Missing, [Type]::Missing, "e2e\e12cf65d-4c4c-4743-92ec-438a640722aa.md")

# ===================== Sheet 2: zh-cn =====================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A5").Value = "d9605128-76f0-495d-8829-707337481735.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "d9605128-76f0-495d-8829-707337481735.9893a002e0b2a800e04aac4f723824fa775d46d4.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-16 20:38:01"

$ws2.Range("A6").Value = "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("G6").Value = "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.eea6052fa4a35afda8e87ade48a954577372908f.zh-cn.xlf"
$ws2.Range("H6").Value = "2016-08-16 20:38:28"

$h2 = $ws2.Hyperlinks
$h2.Delete()
$h2.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b9dc80479d3c1d79b79aaf5b2b4c141c2df2ce0/e2e/4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md", [Type]::Missing, [Type]::Missing, "4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md")
$h2.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/54e1007aeb5a4c815ff1ca79f9d096e0b901819a/e2e/4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md", [Type]::Missing, [Type]::Missing, "4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md")
$h2.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/49145225-09c9-46ae-803f-739b042b9bda.md", [Type]::Missing, [Type]::Missing, "49145225-09c9-46ae-803f-739b042b9bda.md")
$h2.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/b1b21d64-4f4b-4615-91cd-ab12a367e019.md", [Type]::Missing, [Type]::Missing, "b1b21d64-4f4b-4615-91cd-ab12a367e019.md")
$h2.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d6f18f4361b4819f466264699ed58e43a3b7325/e2e/d9605128-76f0-495d-8829-707337481735.md", [Type]::Missing, [Type]::Missing, "d9605128-76f0-495d-8829-707337481735.md")
$h2.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5995be60b2f8ebe30d1c30a9afe1cb90e32440ed/e2e/1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md", [Type]::Missing, [Type]::Missing, "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md")
$h2.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d54b2f816603fc8dda517ffcd058267ebc1ea66/e2e/e12cf65d-4c4c-4743-92ec-438a640722aa.md", [Type]::Missing, [Type]::Missing, "e12cf65d-4c4c-4743-92ec-438a640722aa.md")

# ===================== Sheet 3: de-de =====================
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A5").Value = "d9605128-76f0-495d-8829-707337481735.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "d9605128-76f0-495d-8829-707337481735.9893a002e0b2a800e04aac4f723824fa775d46d4.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-16 20:38:11"

$ws3.Range("A6").Value = "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("G6").Value = "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.eea6052fa4a35afda8e87ade48a954577372908f.de-de.xlf"
$ws3.Range("H6").Value = "2016-08-16 20:38:33"

$h3 = $ws3.Hyperlinks
$h3.Delete()
$h3.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b9dc80479d3c1d79b79aaf5b2b4c141c2df2ce0/e2e/4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md", [Type]::Missing, [Type]::Missing, "4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md")
$h3.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e699b72c65f8deaaabcc3403b8ddaac47a380940/e2e/4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md", [Type]::Missing, [Type]::Missing, "4e594af7-a0fb-4db3-987d-7f24e75b4bb3.md")
$h3.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/49145225-09c9-46ae-803f-739b042b9bda.md", [Type]::Missing, [Type]::Missing, "49145225-09c9-46ae-803f-739b042b9bda.md")
$h3.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84baa329795703b72a6ed1cf19cd66df02c2ff02/e2e/b1b21d64-4f4b-4615-91cd-ab12a367e019.md", [Type]::Missing, [Type]::Missing, "b1b21d64-4f4b-4615-91cd-ab12a367e019.md")
$h3.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d6f18f4361b4819f466264699ed58e43a3b7325/e2e/d9605128-76f0-495d-8829-707337481735.md", [Type]::Missing, [Type]::Missing, "d9605128-76f0-495d-8829-707337481735.md")
$h3.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5995be60b2f8ebe30d1c30a9afe1cb90e32440ed/e2e/1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md", [Type]::Missing, [Type]::Missing, "1d0d76ba-708e-4fcd-a4f2-1bb78e82463c.md")
$h3.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d54b2f816603fc8dda517ffcd058267ebc1ea66/e2e/e12cf65d-4c4c-4743-92ec-438a640722aa.md", [Type]::Missing, [Type]::Missing, "e12cf65d-4c4c-4743-92ec-438a640722aa.md")
